$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (semantic column renames) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case the Spanish connector words ("de"->"De", "del"->"Del", etc.) ---
# --- across the municipality / state name columns ---
$ws.Range('B28').Value = 'Comitán De Domínguez'
$ws.Range('B36').Value = 'Mazapa De Madero'
$ws.Range('B56').Value = 'Coyame Del Sotol'
$ws.Range('B60').Value = 'Hidalgo Del Parral'
$ws.Range('B76').Value = 'Villa De Álvarez'
$ws.Range('A78').Value = 'Ciudad De México'
$ws.Range('B82').Value = 'Cuajimalpa De Morelos'
$ws.Range('B106').Value = 'San Juan Del Río'
$ws.Range('A113').Value = 'Estado De México'
$ws.Range('B113').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B114').Value = 'Almoloya De Juárez'
$ws.Range('B121').Value = 'Chapa De Mota'
$ws.Range('B127').Value = 'Ecatepec De Morelos'
$ws.Range('B129').Value = 'Ixtapan De La Sal'
$ws.Range('B132').Value = 'Naucalpan De Juárez'
$ws.Range('B138').Value = 'San Felipe Del Progreso'
$ws.Range('B145').Value = 'Tenango Del Valle'
$ws.Range('B153').Value = 'Tlalnepantla De Baz'
$ws.Range('B159').Value = 'Valle De Bravo'
$ws.Range('B167').Value = 'Apaseo El Alto'
$ws.Range('B168').Value = 'Apaseo El Grande'
$ws.Range('B174').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B178').Value = 'Jaral Del Progreso'
$ws.Range('B185').Value = 'Purísima Del Rincón'
$ws.Range('B189').Value = 'San Diego De La Unión'
$ws.Range('B191').Value = 'San Luis De La Paz'
$ws.Range('B192').Value = 'Silao De La Victoria'
$ws.Range('B195').Value = 'Valle De Santiago'
$ws.Range('B198').Value = 'Acapulco De Juárez'
$ws.Range('B200').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B201').Value = 'Alcozauca De Guerrero'
$ws.Range('B205').Value = 'Atoyac De Álvarez'
$ws.Range('B206').Value = 'Ayutla De Los Libres'
$ws.Range('B209').Value = 'Buenavista De Cuéllar'
$ws.Range('B210').Value = 'Chilapa De Álvarez'
$ws.Range('B211').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B212').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B217').Value = 'Coyuca De Benítez'
$ws.Range('B218').Value = 'Coyuca De Catalán'
$ws.Range('B221').Value = 'Cuetzala Del Progreso'
$ws.Range('B222').Value = 'Cutzamala De Pinzón'
$ws.Range('B226').Value = 'Iguala De La Independencia'
$ws.Range('B228').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B229').Value = 'Zihuatanejo De Azueta'
$ws.Range('B231').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B234').Value = 'Mártir De Cuilapan'
$ws.Range('B243').Value = 'Taxco De Alarcón'
$ws.Range('B245').Value = 'Técpan De Galeana'
$ws.Range('B247').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B252').Value = 'Tlapa De Comonfort'
$ws.Range('B264').Value = 'Atotonilco El Grande'
$ws.Range('B267').Value = 'Cuautepec De Hinojosa'
$ws.Range('B269').Value = 'Huasca De Ocampo'
$ws.Range('B272').Value = 'Jacala De Ledezma'
$ws.Range('B276').Value = 'Molango De Escamilla'
$ws.Range('B277').Value = 'Omitlán De Juárez'
$ws.Range('B278').Value = 'Pachuca De Soto'
$ws.Range('B284').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B286').Value = 'Tulancingo De Bravo'
$ws.Range('B289').Value = 'Zacualtipán De Ángeles'
$ws.Range('B293').Value = 'Acatlán De Juárez'
$ws.Range('B296').Value = 'Atemajac De Brizuela'
$ws.Range('B298').Value = 'Atotonilco El Alto'
$ws.Range('B299').Value = 'Autlán De Navarro'
$ws.Range('B304').Value = 'Cañadas De Obregón'
$ws.Range('B310').Value = 'Cuautitlán De García Barragán'
$ws.Range('B314').Value = 'Encarnación De Díaz'
$ws.Range('B319').Value = 'Huejuquilla El Alto'
$ws.Range('B320').Value = 'Ixtlahuacán Del Río'
$ws.Range('B323').Value = 'Jilotlán De Los Dolores'
$ws.Range('B326').Value = 'Lagos De Moreno'
$ws.Range('B335').Value = 'San Juan De Los Lagos'
$ws.Range('B337').Value = 'San Martín De Bolaños'
$ws.Range('B339').Value = 'San Miguel El Alto'
$ws.Range('B340').Value = 'Santa María De Los Ángeles'
$ws.Range('B343').Value = 'Tamazula De Gordiano'
$ws.Range('B347').Value = 'Teocuitatlán De Corona'
$ws.Range('B348').Value = 'Tepatitlán De Morelos'
$ws.Range('B356').Value = 'Unión De Tula'
$ws.Range('B359').Value = 'Zacoalco De Torres'
$ws.Range('B361').Value = 'Zapotitlán De Vadillo'
$ws.Range('B362').Value = 'Zapotlán Del Rey'
$ws.Range('B363').Value = 'Zapotlán El Grande'
$ws.Range('B385').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B387').Value = 'Cojumatlán De Régules'
$ws.Range('B447').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B474').Value = 'Jonacatepec De Leandro Valle'
$ws.Range('B477').Value = 'Puente De Ixtla'
$ws.Range('B481').Value = 'Tetela Del Volcán'
$ws.Range('B490').Value = 'Amatlán De Cañas'
$ws.Range('B494').Value = 'Ixtlán Del Río'
$ws.Range('B501').Value = 'Santa María Del Oro'
$ws.Range('B511').Value = 'Montemorelos'
$ws.Range('B519').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B521').Value = 'Coicoyán De Las Flores'
$ws.Range('B522').Value = 'Constancia Del Rosario'
$ws.Range('B525').Value = 'Guadalupe De Ramírez'
$ws.Range('B526').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B527').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B528').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B530').Value = 'Ixtlán De Juárez'
$ws.Range('B531').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B538').Value = 'Mariscala De Juárez'
$ws.Range('B539').Value = 'Mártires De Tacubaya'
$ws.Range('B542').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B543').Value = 'Mixistlán De La Reforma'
$ws.Range('B547').Value = 'Nejapa De Madero'
$ws.Range('B548').Value = 'Oaxaca De Juárez'
$ws.Range('B549').Value = 'Ocotlán De Morelos'
$ws.Range('B550').Value = 'Putla Villa De Guerrero'
$ws.Range('B571').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B585').Value = 'San Juan Bautista Lo De Soto'
$ws.Range('B614').Value = 'San Miguel El Grande'
$ws.Range('B627').Value = 'San Pedro El Alto'
$ws.Range('B635').Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range('B653').Value = 'Santa Inés Del Monte'
$ws.Range('B655').Value = 'Santa Lucía Del Camino'
$ws.Range('B659').Value = 'Santa María Del Tule'
$ws.Range('B662').Value = 'Santa María Jalapa Del Marqués'
$ws.Range('B696').Value = 'Santo Domingo De Morelos'
$ws.Range('B707').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B708').Value = 'Tataltepec De Valdés'
$ws.Range('B709').Value = 'Teotitlán Del Valle'
$ws.Range('B710').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B711').Value = 'Tlacolula De Matamoros'
$ws.Range('B712').Value = 'Tlalixtac De Cabrera'
$ws.Range('B713').Value = 'Totontepec Villa De Morelos'
$ws.Range('B716').Value = 'Villa De Etla'
$ws.Range('B717').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B718').Value = 'Villa De Zaachila'
$ws.Range('B720').Value = 'Zimatlán De Álvarez'
$ws.Range('B733').Value = 'Chalchicomula De Sesma'
$ws.Range('B742').Value = 'Cuayuca De Andrade'
$ws.Range('B743').Value = 'Cuetzalan Del Progreso'
$ws.Range('B753').Value = 'Izúcar De Matamoros'
$ws.Range('B759').Value = 'Los Reyes De Juárez'
$ws.Range('B764').Value = 'Palmar De Bravo'
$ws.Range('B776').Value = 'San Salvador El Verde'
$ws.Range('B784').Value = 'Tepanco De López'
$ws.Range('B785').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B788').Value = 'Tepexi De Rodríguez'
$ws.Range('B789').Value = 'Tetela De Ocampo'
$ws.Range('B793').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B812').Value = 'Amealco De Bonfil'
$ws.Range('B813').Value = 'Cadereyta De Montes'
$ws.Range('B818').Value = 'Pinal De Amoles'
$ws.Range('B820').Value = 'San Juan Del Río'
$ws.Range('B830').Value = 'Mexquitic De Carmona'
$ws.Range('B833').Value = 'San Ciro De Acosta'
$ws.Range('B835').Value = 'Santa María Del Río'
$ws.Range('B839').Value = 'Villa De Ramos'
$ws.Range('B903').Value = 'Amatlán De Los Reyes'
$ws.Range('B909').Value = 'Boca Del Río'
$ws.Range('B914').Value = 'Castillo De Teayo'
$ws.Range('B922').Value = 'Cosamaloapan De Carpio'
$ws.Range('B934').Value = 'Ixhuatlán De Madero'
$ws.Range('B939').Value = 'Juchique De Ferrer'
$ws.Range('B945').Value = 'Martínez De La Torre'
$ws.Range('B946').Value = 'Medellín De Bravo'
$ws.Range('B949').Value = 'Mixtla De Altamirano'
$ws.Range('B960').Value = 'Sayula De Alemán'
$ws.Range('B963').Value = 'Soledad De Doblado'
$ws.Range('B1008').Value = 'Concepción Del Oro'
$ws.Range('B1018').Value = 'Nochistlán De Mejía'
$ws.Range('B1026').Value = 'Teúl De González Ortega'
$ws.Range('B1027').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1028').Value = 'Trinidad García De La Cadena'
$ws.Range('B1030').Value = 'Villa De Cos'

# --- Drop the trailing footnote / metadata rows (1037-1041) ---
$ws.Range('A1037:A1041').ClearContents()
